$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values from re-pulled data
$ws.Range("F2").Value = -9
$ws.Range("F4").Value = -2
$ws.Range("F6").Value = -4
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 2
$ws.Range("F12").Value = -3
$ws.Range("F13").Value = -7
$ws.Range("F14").Value = 4
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = -1
$ws.Range("F19").Value = -1
$ws.Range("F20").Value = 1
$ws.Range("F22").Value = -3
$ws.Range("F23").Value = -3
